$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.3465496666666667
$ws.Cells.Item(2, 8).Value = 1.039649
$ws.Cells.Item(2, 9).Value = 0.008996151488293185
$ws.Cells.Item(2, 10).Value = 0.008996151488293185
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.3465496666666667
$ws.Cells.Item(2, 14).Value = 1.039649
$ws.Cells.Item(2, 15).Value = 0.008996151488293185
$ws.Cells.Item(2, 16).Value = 0.008996151488293185
$ws.Cells.Item(2, 17).Value = 0.1200966714667778
$ws.Cells.Item(2, 18).Value = 1.080870043201
$ws.Cells.Item(2, 19).Value = 0.00008093074160031967
$ws.Cells.Item(2, 20).Value = 0.00008093074160031967

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.3465496666666667
$ws.Cells.Item(3, 8).Value = 1.039649
$ws.Cells.Item(3, 9).Value = 0.008996151488293185
$ws.Cells.Item(3, 10).Value = 0.008996151488293185
$ws.Cells.Item(3, 14).Value = 66.23320799999999
$ws.Cells.Item(3, 15).Value = 0.5731203249593199
$ws.Cells.Item(3, 16).Value = 0.5731203249593199
$ws.Cells.Item(3, 17).Value = 7.651032051554667
$ws.Cells.Item(3, 18).Value = 68.85928846399199
$ws.Cells.Item(3, 19).Value = 0.005155877264353859
$ws.Cells.Item(3, 20).Value = 0.005155877264353859

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.3465496666666667
$ws.Cells.Item(4, 8).Value = 1.039649
$ws.Cells.Item(4, 9).Value = 0.008996151488293185
$ws.Cells.Item(4, 10).Value = 0.008996151488293185
$ws.Cells.Item(4, 13).Value = 16.08941833333333
$ws.Cells.Item(4, 14).Value = 48.268255
$ws.Cells.Item(4, 15).Value = 0.4176683996767803
$ws.Cells.Item(4, 16).Value = 0.4176683996767803
$ws.Cells.Item(4, 17).Value = 5.575782560277223
$ws.Cells.Item(4, 18).Value = 50.182043042495
$ws.Cells.Item(4, 19).Value = 0.003757408195365299
$ws.Cells.Item(4, 20).Value = 0.0037574081953653

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.3465496666666667
$ws.Cells.Item(5, 8).Value = 1.039649
$ws.Cells.Item(5, 9).Value = 0.008996151488293185
$ws.Cells.Item(5, 10).Value = 0.008996151488293185
$ws.Cells.Item(5, 13).Value = 0.008287000000000001
$ws.Cells.Item(5, 14).Value = 0.024861
$ws.Cells.Item(5, 15).Value = 0.0002151238756065334
$ws.Cells.Item(5, 16).Value = 0.0002151238756065334
$ws.Cells.Item(5, 17).Value = 0.002871857087666667
$ws.Cells.Item(5, 18).Value = 0.025846713789
$ws.Cells.Item(5, 19).Value = 0.000001935286973705113
$ws.Cells.Item(5, 20).Value = 0.000001935286973705114

$ws.Cells.Item(6, 8).Value = 66.23320799999999
$ws.Cells.Item(6, 9).Value = 0.5731203249593199
$ws.Cells.Item(6, 10).Value = 0.5731203249593199
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.3465496666666667
$ws.Cells.Item(6, 14).Value = 1.039649
$ws.Cells.Item(6, 15).Value = 0.008996151488293185
$ws.Cells.Item(6, 16).Value = 0.008996151488293185
$ws.Cells.Item(6, 17).Value = 7.651032051554667
$ws.Cells.Item(6, 18).Value = 68.85928846399199
$ws.Cells.Item(6, 19).Value = 0.005155877264353859
$ws.Cells.Item(6, 20).Value = 0.005155877264353859

$ws.Cells.Item(7, 8).Value = 66.23320799999999
$ws.Cells.Item(7, 9).Value = 0.5731203249593199
$ws.Cells.Item(7, 10).Value = 0.5731203249593199
$ws.Cells.Item(7, 14).Value = 66.23320799999999
$ws.Cells.Item(7, 15).Value = 0.5731203249593199
$ws.Cells.Item(7, 16).Value = 0.5731203249593199
$ws.Cells.Item(7, 17).Value = 487.4264268856959
$ws.Cells.Item(7, 18).Value = 4386.837841971263
$ws.Cells.Item(7, 19).Value = 0.3284669068814764
$ws.Cells.Item(7, 20).Value = 0.3284669068814764

$ws.Cells.Item(8, 8).Value = 66.23320799999999
$ws.Cells.Item(8, 9).Value = 0.5731203249593199
$ws.Cells.Item(8, 10).Value = 0.5731203249593199
$ws.Cells.Item(8, 13).Value = 16.08941833333333
$ws.Cells.Item(8, 14).Value = 48.268255
$ws.Cells.Item(8, 15).Value = 0.4176683996767803
$ws.Cells.Item(8, 16).Value = 0.4176683996767803
$ws.Cells.Item(8, 17).Value = 355.2179303568933
$ws.Cells.Item(8, 18).Value = 3196.96137321204
$ws.Cells.Item(8, 19).Value = 0.2393742489479954
$ws.Cells.Item(8, 20).Value = 0.2393742489479954

$ws.Cells.Item(9, 8).Value = 66.23320799999999
$ws.Cells.Item(9, 9).Value = 0.5731203249593199
$ws.Cells.Item(9, 10).Value = 0.5731203249593199
$ws.Cells.Item(9, 13).Value = 0.008287000000000001
$ws.Cells.Item(9, 14).Value = 0.024861
$ws.Cells.Item(9, 15).Value = 0.0002151238756065334
$ws.Cells.Item(9, 16).Value = 0.0002151238756065334
$ws.Cells.Item(9, 17).Value = 0.182958198232
$ws.Cells.Item(9, 18).Value = 1.646623784088
$ws.Cells.Item(9, 19).Value = 0.0001232918654941247
$ws.Cells.Item(9, 20).Value = 0.0001232918654941248

$ws.Cells.Item(10, 7).Value = 16.08941833333333
$ws.Cells.Item(10, 8).Value = 48.268255
$ws.Cells.Item(10, 9).Value = 0.4176683996767803
$ws.Cells.Item(10, 10).Value = 0.4176683996767803
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.3465496666666667
$ws.Cells.Item(10, 14).Value = 1.039649
$ws.Cells.Item(10, 15).Value = 0.008996151488293185
$ws.Cells.Item(10, 16).Value = 0.008996151488293185
$ws.Cells.Item(10, 17).Value = 5.575782560277223
$ws.Cells.Item(10, 18).Value = 50.182043042495
$ws.Cells.Item(10, 19).Value = 0.003757408195365299
$ws.Cells.Item(10, 20).Value = 0.0037574081953653

$ws.Cells.Item(11, 7).Value = 16.08941833333333
$ws.Cells.Item(11, 8).Value = 48.268255
$ws.Cells.Item(11, 9).Value = 0.4176683996767803
$ws.Cells.Item(11, 10).Value = 0.4176683996767803
$ws.Cells.Item(11, 14).Value = 66.23320799999999
$ws.Cells.Item(11, 15).Value = 0.5731203249593199
$ws.Cells.Item(11, 16).Value = 0.5731203249593199
$ws.Cells.Item(11, 17).Value = 355.2179303568933
$ws.Cells.Item(11, 18).Value = 3196.96137321204
$ws.Cells.Item(11, 19).Value = 0.2393742489479954
$ws.Cells.Item(11, 20).Value = 0.2393742489479954

$ws.Cells.Item(12, 7).Value = 16.08941833333333
$ws.Cells.Item(12, 8).Value = 48.268255
$ws.Cells.Item(12, 9).Value = 0.4176683996767803
$ws.Cells.Item(12, 10).Value = 0.4176683996767803
$ws.Cells.Item(12, 13).Value = 16.08941833333333
$ws.Cells.Item(12, 14).Value = 48.268255
$ws.Cells.Item(12, 15).Value = 0.4176683996767803
$ws.Cells.Item(12, 16).Value = 0.4176683996767803
$ws.Cells.Item(12, 17).Value = 258.8693823050028
$ws.Cells.Item(12, 18).Value = 2329.824440745025
$ws.Cells.Item(12, 19).Value = 0.1744468920885626
$ws.Cells.Item(12, 20).Value = 0.1744468920885627

$ws.Cells.Item(13, 7).Value = 16.08941833333333
$ws.Cells.Item(13, 8).Value = 48.268255
$ws.Cells.Item(13, 9).Value = 0.4176683996767803
$ws.Cells.Item(13, 10).Value = 0.4176683996767803
$ws.Cells.Item(13, 13).Value = 0.008287000000000001
$ws.Cells.Item(13, 14).Value = 0.024861
$ws.Cells.Item(13, 15).Value = 0.0002151238756065334
$ws.Cells.Item(13, 16).Value = 0.0002151238756065334
$ws.Cells.Item(13, 17).Value = 0.1333330097283334
$ws.Cells.Item(13, 18).Value = 1.199997087555
$ws.Cells.Item(13, 19).Value = 0.00008985044485684756
$ws.Cells.Item(13, 20).Value = 0.00008985044485684758

$ws.Cells.Item(14, 7).Value = 0.008287000000000001
$ws.Cells.Item(14, 8).Value = 0.024861
$ws.Cells.Item(14, 9).Value = 0.0002151238756065334
$ws.Cells.Item(14, 10).Value = 0.0002151238756065334
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.3465496666666667
$ws.Cells.Item(14, 14).Value = 1.039649
$ws.Cells.Item(14, 15).Value = 0.008996151488293185
$ws.Cells.Item(14, 16).Value = 0.008996151488293185
$ws.Cells.Item(14, 17).Value = 0.002871857087666667
$ws.Cells.Item(14, 18).Value = 0.025846713789
$ws.Cells.Item(14, 19).Value = 0.000001935286973705113
$ws.Cells.Item(14, 20).Value = 0.000001935286973705114

$ws.Cells.Item(15, 7).Value = 0.008287000000000001
$ws.Cells.Item(15, 8).Value = 0.024861
$ws.Cells.Item(15, 9).Value = 0.0002151238756065334
$ws.Cells.Item(15, 10).Value = 0.0002151238756065334
$ws.Cells.Item(15, 14).Value = 66.23320799999999
$ws.Cells.Item(15, 15).Value = 0.5731203249593199
$ws.Cells.Item(15, 16).Value = 0.5731203249593199
$ws.Cells.Item(15, 17).Value = 0.182958198232
$ws.Cells.Item(15, 18).Value = 1.646623784088
$ws.Cells.Item(15, 19).Value = 0.0001232918654941247
$ws.Cells.Item(15, 20).Value = 0.0001232918654941248

$ws.Cells.Item(16, 7).Value = 0.008287000000000001
$ws.Cells.Item(16, 8).Value = 0.024861
$ws.Cells.Item(16, 9).Value = 0.0002151238756065334
$ws.Cells.Item(16, 10).Value = 0.0002151238756065334
$ws.Cells.Item(16, 13).Value = 16.08941833333333
$ws.Cells.Item(16, 14).Value = 48.268255
$ws.Cells.Item(16, 15).Value = 0.4176683996767803
$ws.Cells.Item(16, 16).Value = 0.4176683996767803
$ws.Cells.Item(16, 17).Value = 0.1333330097283334
$ws.Cells.Item(16, 18).Value = 1.199997087555
$ws.Cells.Item(16, 19).Value = 0.00008985044485684756
$ws.Cells.Item(16, 20).Value = 0.00008985044485684758

$ws.Cells.Item(17, 7).Value = 0.008287000000000001
$ws.Cells.Item(17, 8).Value = 0.024861
$ws.Cells.Item(17, 9).Value = 0.0002151238756065334
$ws.Cells.Item(17, 10).Value = 0.0002151238756065334
$ws.Cells.Item(17, 13).Value = 0.008287000000000001
$ws.Cells.Item(17, 14).Value = 0.024861
$ws.Cells.Item(17, 15).Value = 0.0002151238756065334
$ws.Cells.Item(17, 16).Value = 0.0002151238756065334
$ws.Cells.Item(17, 17).Value = 0.00006867436900000002
$ws.Cells.Item(17, 18).Value = 0.0006180693210000001
$ws.Cells.Item(17, 19).Value = 0.00000004627828185597526
$ws.Cells.Item(17, 20).Value = 0.00000004627828185597527
